$wb = $excel.ActiveWorkbook

# Mapping of sheet name -> row number -> new F-column value
$sheet1Updates = @{
    2  = 275
    3  = 576
    5  = 288
    6  = 1109
    7  = 1445
    8  = 588
    12 = 171
    13 = 121
    14 = 445
    15 = 1382
    16 = 119
    17 = 118
    20 = 57
    21 = 659
    24 = 233
    25 = 23
    26 = 5941
    28 = 121
    29 = 111
    31 = 14615
    32 = 1448
    33 = 220
    36 = 9176
    37 = 634
    38 = 4219
    39 = 152
}

$sheet4Updates = @{
    2  = 275
    3  = 576
    5  = 288
    6  = 1109
    7  = 1445
    8  = 588
    12 = 171
    13 = 121
    14 = 445
    15 = 1382
    16 = 119
    17 = 118
    21 = 57
    22 = 659
    26 = 233
    27 = 23
    29 = 5941
    31 = 121
    32 = 111
    34 = 14615
    35 = 1448
    36 = 220
    39 = 9176
    40 = 634
    41 = 4219
    42 = 152
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
